$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text change: "Ready for handoff" -> "In Translation"
#    This shared string is referenced from the Overview sheet (columns for
#    zh-cn / de-de status) as well as the per-locale detail sheets, so every
#    cell holding that value needs to be updated.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the status columns (Overview!E:F and the "Status" column on each
#    locale sheet) from ~17.22 chars down to ~13.41 chars.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
